# Updated symbol list on Sun Dec 25 11:18:02 UTC 2022 with GitHub Actions
#
# The "Price" column (D) stores numbers formatted as plain text (e.g. "244.80",
# "0.0009838") so the original string representation -- including trailing
# zeros -- must be preserved exactly. Assigning a numeric-looking string to
# a General-formatted cell would silently convert it to a real number and
# drop significant trailing zeros, so we prefix the literal with a leading
# apostrophe (exactly like a user typing a text-forced value into Excel) to
# keep these cells as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $ws.Range($range).Value = "'" + $value
}

# --- price refreshes (Price column, D) ---
Set-TextValue "D2"  "244.78"
Set-TextValue "D4"  "5.432"
Set-TextValue "D5"  "0.05971"
Set-TextValue "D7"  "0.8095"
Set-TextValue "D8"  "0.9262"
Set-TextValue "D9"  "0.1431"
Set-TextValue "D10" "0.07442"
Set-TextValue "D11" "0.03369"
Set-TextValue "D12" "0.03036"
Set-TextValue "D13" "0.09346"
Set-TextValue "D14" "3.946"
Set-TextValue "D15" "0.001593"
Set-TextValue "D16" "0.04808"
Set-TextValue "D17" "0.0005944"
Set-TextValue "D18" "0.005698"
Set-TextValue "D19" "0.004155"
Set-TextValue "D20" "0.0009820"
Set-TextValue "D21" "0.00007704"
Set-TextValue "D22" "3.659"
Set-TextValue "D23" "6.451"
Set-TextValue "D26" "0.1337"
Set-TextValue "D27" "0.0002448"
Set-TextValue "D40" "0.03939"

# --- rows 41-43: symbol list rotated (KickToken/BKEXToken/CEJI shifted up) ---
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.006213"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1073"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.002681"
$ws.Range("E43").Value = "42CEJICEJI"

# --- remaining price refreshes ---
Set-TextValue "D44" "0.007244"
Set-TextValue "D45" "0.00005130"
Set-TextValue "D47" "0.0005804"
Set-TextValue "D49" "0.002279"
